$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Analisar prestação de contas" test case moves from the TC4 slot to the
# TC3 slot (3rd block, rows 14-18), pushing "Detalhar diária" (was TC2 slot)
# down to the TC3 slot (rows 21-25) and "Cancelar diária" (was TC3 slot) down
# to the TC4 slot (rows 28-32). The "TCx" id labels themselves (B14, B21,
# B28) stay fixed in place; only the step/expected-result text moves.

$ws.Range("B18").Value = "Beneficiário Clica em analisar prestação de contas."
$ws.Range("D18").Value = "SYSTEM Apresenta a tela de Analisar Prestação de Contas"

$ws.Range("B25").Value = "Beneficiário Clica em detalhar diária."
$ws.Range("D25").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"

$ws.Range("B32").Value = "Beneficiário Clica em cancelar diária."
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária"
